$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 45
$ws.Range("I2").Value = 124
$ws.Range("J2").Value = 488
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 131
$ws.Range("M2").Value = 11
$ws.Range("N2").Value = 70
$ws.Range("P2").Value = 4
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 6
$ws.Range("S2").Value = 58
$ws.Range("T2").Value = 91
$ws.Range("U2").Value = 7
$ws.Range("V2").Value = 705
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 757
$ws.Range("Y2").Value = 3
$ws.Range("Z2").Value = 6
